$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: average of the "|S*|/n" column (J) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Row 14: Average of SW(S*)/SW(OPT) ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$f14 = $ws.Range("B14").Font
$f14.Bold = $true
$f14.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# --- Row 15: Average of SC(S*)/SC(OPT) ---
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$f15 = $ws.Range("B15").Font
$f15.Bold = $true
$f15.Size = 12
$ws.Range("B15").VerticalAlignment = -4108

# --- Row 16: Worst of SW(S*)/SW(OPT) ---
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$f16 = $ws.Range("B16").Font
$f16.Bold = $true
$f16.Size = 12
$ws.Range("B16").VerticalAlignment = -4108

# --- Row 17: Worst of SC(S*)/SC(OPT) ---
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$f17 = $ws.Range("B17").Font
$f17.Bold = $true
$f17.Size = 12
$ws.Range("B17").VerticalAlignment = -4108

# Taller rows for the summary block (matches ht="15.6" in the authored file)
$ws.Range("A14:B17").EntireRow.RowHeight = 15.6

# --- Selection cursor left on J12 by the editing session ---
$ws.Range("J12").Select()

# --- Page setup (paper size / orientation) added when saved from this session ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "edit complete"
